$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Id"
$ws.Range("C1").Value = "Task Name"

$ws.Range("C2").Select()
